$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tables")

# Require schema to be explicitly given in the DB table definitions:
# prefix the "ref_*" / "constants" lookup column (B) with "reference." so
# it reads as "reference.ref_xxx" / "reference.constants".
$ws.Range("B2").Value = "reference.ref_general_parameter"
$ws.Range("B3").Value = "reference.ref_drag_coef_cyl"
$ws.Range("B4").Value = "reference.ref_wake_amplification_factor_cyl"
$ws.Range("B5").Value = "reference.ref_wind_drag_coef_rect"
$ws.Range("B6").Value = "reference.ref_current_drag_coef_rect"
$ws.Range("B7").Value = "reference.ref_drift_coef_float_rect"
$ws.Range("B8").Value = "reference.ref_rectangular_wave_inertia"
$ws.Range("B9").Value = "reference.constants"
$ws.Range("B10").Value = "reference.constants"
$ws.Range("B11").Value = "reference.constants"
$ws.Range("B12").Value = "reference.constants"
$ws.Range("B13").Value = "reference.constants"
$ws.Range("B14").Value = "reference.constants"

# Widen column B to fit the longer "reference.*" values.
$ws.Columns.Item(2).ColumnWidth = 42.7

# Give row 4 an explicit height (as in the edited workbook).
$ws.Rows.Item(4).RowHeight = 15.5

# Make "Tables" the active sheet/tab, with B14 selected.
$ws.Activate()
$ws.Range("B14").Select()
